# Adds three new RFLTools commands to the command-reference worksheet:
#   GETZ2 (inserted alphabetically between GETGRID and GEXSUPER, i.e. at row 60)
#   Z2N and Z2P (appended after the last row, Z2ATT, at the end of the list)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert GETZ2 just above the current row 60 (GEXSUPER), pushing
#     GEXSUPER and everything below it down by one row. ---
$ws.Rows("60:60").Insert()
$ws.Range("A60").Value = "GETZ2"
$ws.Range("B60").Value = "Searches a block or xref for nested blocks (useful for finding all blocks such as signs).  Used with Z2N and Z2P"

# --- Append Z2N and Z2P after the (now shifted) last row, Z2ATT at row 184. ---
$ws.Range("A185").Value = "Z2N"
$ws.Range("B185").Value = "Zooms to the next block found with GETZ2"
$ws.Range("A186").Value = "Z2P"
$ws.Range("B186").Value = "Zooms to the previous block found with GETZ2"

# --- Restore the view state (scroll position / active cell) to match the
#     saved workbook's final sheetView. ---
$ws.Range("A128").Select()
$ws.Range("B187").Select()
